$d = $word.ActiveDocument
$sel = $word.Selection

# Jump to the end of the document body (after "...my test document")
$sel.EndKey(6) | Out-Null  # wdStory

# Start a brand-new paragraph and make sure it carries no inherited direct
# character formatting (e.g. the sz/szCs that sits on the first paragraph's
# mark) before applying the Heading 1 style.
$sel.TypeParagraph() | Out-Null
$sel.ClearFormatting() | Out-Null
$sel.Style = "Heading1"
$sel.TypeText("Test Heading with Violations")

# Format just the heading's text run (Arial + the MACE navy accent color),
# leaving the paragraph mark's run properties untouched.
$newPara = $d.Paragraphs.Last
$fullRange = $newPara.Range
$textRange = $d.Range($fullRange.Start, $fullRange.End - 1)
$textRange.Font.Name = "Arial"
$textRange.Font.Color = 0x993300  # wdColor (BGR) for RGB 003399
